$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry for 27/12/2024: a day of "UI" work, same row-formatting
# pattern as the existing entries (copy row 21's layout down onto the
# freshly appended row 22, then overwrite with the real data).
$ws.Range("A21:E21").Copy()
$ws.Range("A22").PasteSpecial()

$ws.Range("A22").Value = 45653
$ws.Range("B22").Value = "UI"
$ws.Range("C22").Value = "Icone, problemi vari"
$ws.Range("D22").Value = 0.083333333333333329
$ws.Range("E22").Value = "Icone mentore, icona prosegui testo, correzioni varie"
